# Auto-generated Excel COM-interop script to apply the cryptos.xlsx update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so values like '17.80' or '46.16'
# are not silently reinterpreted/rounded as numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '26.351.88'
$ws.Range("E2").Value = '  -1.94%  '

# Row 3
$ws.Range("D3").Value = '1.792.86'
$ws.Range("E3").Value = '  -1.79%  '

# Row 4
$ws.Range("D4").Value = '1.005'
$ws.Range("E4").Value = '  -0.20%  '

# Row 5
$ws.Range("D5").Value = '1.005'
$ws.Range("E5").Value = '  -0.16%  '

# Row 6
$ws.Range("D6").Value = '306.96'
$ws.Range("E6").Value = '  -1.27%  '

# Row 7
$ws.Range("D7").Value = '0.4522'
$ws.Range("E7").Value = '  -1.21%  '

# Row 8
$ws.Range("D8").Value = '0.3590'
$ws.Range("E8").Value = '  -2.30%  '

# Row 9
$ws.Range("B9").Value = 'OKB'
$ws.Range("C9").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D9").Value = '46.16'
$ws.Range("E9").Value = '  +0.72%  '

# Row 10
$ws.Range("B10").Value = 'Dogecoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D10").Value = '0.07071'
$ws.Range("E10").Value = '  -1.18%  '

# Row 11
$ws.Range("B11").Value = 'Polygon'
$ws.Range("C11").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D11").Value = '0.8840'
$ws.Range("E11").Value = '  +1.40%  '

# Row 12
$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").Value = '0.07762'
$ws.Range("E12").Value = '  -0.32%  '

# Row 13
$ws.Range("B13").Value = 'Solana'
$ws.Range("C13").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D13").Value = '19.44'
$ws.Range("E13").Value = '  -0.44%  '

# Row 14
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.818.00'
$ws.Range("E14").Value = '  -0.05%  '

# Row 15
$ws.Range("B15").Value = 'Polkadot'
$ws.Range("C15").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D15").Value = '5.281'
$ws.Range("E15").Value = '  -0.68%  '

# Row 16
$ws.Range("B16").Value = 'Chainlink'
$ws.Range("C16").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D16").Value = '6.317'
$ws.Range("E16").Value = '  -0.97%  '

# Row 17
$ws.Range("B17").Value = 'Litecoin'
$ws.Range("C17").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D17").Value = '84.96'
$ws.Range("E17").Value = '  -2.13%  '

# Row 18
$ws.Range("B18").Value = 'BinanceUSD'
$ws.Range("C18").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D18").Value = '1.006'
$ws.Range("E18").Value = '  -0.20%  '

# Row 19
$ws.Range("B19").Value = 'ShibaInu'
$ws.Range("C19").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D19").Value = '0.000008501'
$ws.Range("E19").Value = '  -2.21%  '

# Row 20
$ws.Range("B20").Value = 'Dai'
$ws.Range("C20").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D20").Value = '1.004'
$ws.Range("E20").Value = '  -0.21%  '

# Row 21
$ws.Range("B21").Value = 'WrappedBTC'
$ws.Range("C21").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D21").Value = '26.379.43'
$ws.Range("E21").Value = '  -1.90%  '

# Row 22
$ws.Range("B22").Value = 'Avalanche'
$ws.Range("C22").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D22").Value = '14.23'
$ws.Range("E22").Value = '  -1.46%  '

# Row 23
$ws.Range("B23").Value = 'Uniswap'
$ws.Range("C23").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D23").Value = '4.964'
$ws.Range("E23").Value = '  -0.54%  '

# Row 24
$ws.Range("D24").Value = '2.032.48'
$ws.Range("E24").Value = '  -0.20%  '

# Row 25
$ws.Range("B25").Value = 'Cosmos'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D25").Value = '10.53'
$ws.Range("E25").Value = '  +0.80%  '

# Row 26
$ws.Range("B26").Value = 'Toncoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D26").Value = '1.966'
$ws.Range("E26").Value = '  -2.12%  '

# Row 27
$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").Value = '151.25'
$ws.Range("E27").Value = '  +0.24%  '

# Row 28
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").Value = '17.80'
$ws.Range("E28").Value = '  -2.00%  '

# Row 29
$ws.Range("B29").Value = 'LidoDAOToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D29").Value = '2.015'
$ws.Range("E29").Value = '  +3.26%  '

# Row 30
$ws.Range("B30").Value = 'BitcoinCash'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D30").Value = '111.84'
$ws.Range("E30").Value = '  -1.37%  '

# Row 31
$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").Value = '4.839'
$ws.Range("E31").Value = '  -1.36%  '

# Row 32
$ws.Range("B32").Value = 'Stellar'
$ws.Range("C32").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D32").Value = '0.08688'
$ws.Range("E32").Value = '  -1.31%  '

# Row 33
$ws.Range("B33").Value = 'HuobiToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D33").Value = '3.116'
$ws.Range("E33").Value = '  +4.06%  '

# Row 34
$ws.Range("B34").Value = 'RenderToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D34").Value = '2.770'
$ws.Range("E34").Value = '  +9.91%  '

# Row 35
$ws.Range("B35").Value = 'Filecoin'
$ws.Range("C35").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D35").Value = '4.437'
$ws.Range("E35").Value = '  -0.70%  '

# Row 36
$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").Value = '0.7194'
$ws.Range("E36").Value = '  -3.70%  '

# Row 37
$ws.Range("B37").Value = 'ARBITRUM'
$ws.Range("C37").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D37").Value = '1.101'
$ws.Range("E37").Value = '  -2.49%  '

# Row 38
$ws.Range("B38").Value = 'Frax'
$ws.Range("C38").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D38").Value = '1.003'
$ws.Range("E38").Value = '  +0.22%  '

# Row 39
$ws.Range("B39").Value = 'TrustWalletToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D39").Value = '1.066'
$ws.Range("E39").Value = '  -1.62%  '

# Row 40
$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").Value = '0.01928'
$ws.Range("E40").Value = '  -0.33%  '

# Row 41
$ws.Range("B41").Value = 'Hedera'
$ws.Range("C41").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D41").Value = '0.05072'
$ws.Range("E41").Value = '  -0.65%  '

# Row 42
$ws.Range("B42").Value = 'MXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D42").Value = '2.856'
$ws.Range("E42").Value = '  -1.83%  '

# Row 43
$ws.Range("B43").Value = 'TheSandbox'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D43").Value = '0.5057'
$ws.Range("E43").Value = '  +1.95%  '

# Row 44
$ws.Range("B44").Value = 'FraxShare'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D44").Value = '6.825'
$ws.Range("E44").Value = '  -1.30%  '

# Row 45
$ws.Range("B45").Value = 'Algorand'
$ws.Range("C45").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D45").Value = '0.1513'
$ws.Range("E45").Value = '  -4.97%  '

# Row 46
$ws.Range("B46").Value = 'Aptos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D46").Value = '8.010'
$ws.Range("E46").Value = '  -3.03%  '

# Row 47
$ws.Range("B47").Value = 'PaxDollar'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D47").Value = '1.005'
$ws.Range("E47").Value = '  -0.14%  '

# Row 48
$ws.Range("B48").Value = 'Decentraland'
$ws.Range("C48").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D48").Value = '0.4620'
$ws.Range("E48").Value = '  -1.09%  '

# Row 49
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = '9.886'
$ws.Range("E49").Value = '  -2.13%  '

# Row 50
$ws.Range("B50").Value = 'Quant'
$ws.Range("C50").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D50").Value = '100.60'
$ws.Range("E50").Value = '  -0.61%  '

# Row 51
$ws.Range("B51").Value = 'NEARProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D51").Value = '1.563'
$ws.Range("E51").Value = '  -2.57%  '

